# "Drop in RMI script results for 3.0"
# EPA moved its Mortality Risk Valuation FAQ page; update the link on the
# About sheet (cell B6) to point at the new URL, both as a live hyperlink
# and as the cell's displayed text.

$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")

$newUrl = "https://www.epa.gov/environmental-economics/mortality-risk-valuation"
$anchor = "whatvalue"
$newUrlWithAnchor = "$newUrl#$anchor"
$displayText = "$newUrl - $anchor"

$cell = $about.Range("B6")

# Create the hyperlink first (this also touches the cell's displayed
# text/style), then overwrite the cell text with the plain URL + anchor
# so the stored value matches the new FAQ address, and reapply the
# built-in Hyperlink style so the formatting stays the way it was.
$about.Hyperlinks.Add($cell, $newUrl, $anchor, [System.Reflection.Missing]::Value, $displayText) | Out-Null
$cell.Value = $newUrlWithAnchor
$cell.Style = "Hyperlink"

$wb.Save()
